$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helpers
# ---------------------------------------------------------------------------

# Write a value that looks like a date ("dd/mm/yyyy") while keeping it as
# literal text (matching the original shared-string cell type), by staging it
# in a scratch cell formatted as Text, copying, and pasting values-only into
# the destination (destination keeps its own "General" number format/style).
function Set-LiteralText {
    param($rangeAddress, $text)
    $scratch = $ws.Range("ZZ1")
    $scratch.NumberFormat = "@"
    $scratch.Value2 = $text
    $scratch.Copy()
    $ws.Range($rangeAddress).PasteSpecial(-4163)  # xlPasteValues
    $scratch.Clear()
}

# Give a brand-new (previously empty) cell the same column formatting as an
# existing, correctly-styled cell in the same column, then write its value.
function Set-NewCell {
    param($rangeAddress, $sourceAddress, $text)
    $ws.Range($sourceAddress).Copy()
    $ws.Range($rangeAddress).PasteSpecial(-4122)  # xlPasteFormats
    $ws.Range($rangeAddress).Value2 = $text
}

# ---------------------------------------------------------------------------
# Reset the row height of rows 17 and 22 back to the sheet default (no
# customHeight marker) by deleting and re-inserting them in place, *before*
# any new content is written into them.
# ---------------------------------------------------------------------------
$ws.Rows.Item(17).Delete()
$ws.Rows.Item(17).Insert()
$ws.Rows.Item(22).Delete()
$ws.Rows.Item(22).Insert()

# ---------------------------------------------------------------------------
# Row 10: Objetivos value swapped for the docente line
# ---------------------------------------------------------------------------
$ws.Range("B10").Value2 = "7926291 - Célia Regina Tomachuk dos Santos Catuogno"
$ws.Range("C10").Value2 = "7926291 - Célia Regina Tomachuk dos Santos Catuogno"

# ---------------------------------------------------------------------------
# Row 13: gains an "A" label, B/C content becomes "Semestral"
# ---------------------------------------------------------------------------
$ws.Range("A13").Value2 = "Programa resumido:"
$ws.Range("B13").Value2 = "Semestral"
$ws.Range("C13").Value2 = "Semestral"
$ws.Rows.Item(13).RowHeight = 60

# ---------------------------------------------------------------------------
# Row 14: label -> "Short syllabus:", content -> EN objectives text
# ---------------------------------------------------------------------------
$ws.Range("A14").Value2 = "Short syllabus:"
$ws.Range("B14").Value2 = "Project definition. Life cycle of projects. Projectos management processes. Sustainable development. Case studies"
$ws.Range("C14").Value2 = "Project definition. Life cycle of projects. Projectos management processes. Sustainable development. Case studies"

# ---------------------------------------------------------------------------
# Row 15: label -> "Programa:", content -> "01/01/2018" (kept as literal text)
# ---------------------------------------------------------------------------
$ws.Range("A15").Value2 = "Programa:"
Set-LiteralText "B15" "01/01/2018"
Set-LiteralText "C15" "01/01/2018"
$ws.Rows.Item(15).RowHeight = 120

# ---------------------------------------------------------------------------
# Row 16: label -> "Syllabus:", content -> EN long program text
# ---------------------------------------------------------------------------
$ws.Range("A16").Value2 = "Syllabus:"
$ws.Range("B16").Value2 = "Definition of project and its main attributes and characteristics. Project Management Body of Knowledge concepts. Strategic planning. Sustainable Development: Definition. Agreements, treaties and policies international carried out around sustainable development. The challenges of sustainable. Processes and Methodology of management of environmental projects. Planning tools, monitoring and control. Study of the risks and common problems in the management of environmental projects. Analysis of cases involving selection, management and development of environmental management applied projects."
$ws.Range("C16").Value2 = "Definition of project and its main attributes and characteristics. Project Management Body of Knowledge concepts. Strategic planning. Sustainable Development: Definition. Agreements, treaties and policies international carried out around sustainable development. The challenges of sustainable. Processes and Methodology of management of environmental projects. Planning tools, monitoring and control. Study of the risks and common problems in the management of environmental projects. Analysis of cases involving selection, management and development of environmental management applied projects."

# ---------------------------------------------------------------------------
# Row 17: label -> "Avaliação:" only (B/C stay empty, height already reset)
# ---------------------------------------------------------------------------
$ws.Range("A17").Value2 = "Avaliação:"

# ---------------------------------------------------------------------------
# Row 18: label -> "Método:", gains docente line in B/C (new cells)
# ---------------------------------------------------------------------------
$ws.Range("A18").Value2 = "Método:"
Set-NewCell "B18" "B19" "7926291 - Célia Regina Tomachuk dos Santos Catuogno"
Set-NewCell "C18" "C19" "7926291 - Célia Regina Tomachuk dos Santos Catuogno"
$ws.Rows.Item(18).RowHeight = 60

# ---------------------------------------------------------------------------
# Row 19: label -> "Critério:" (B/C unchanged)
# ---------------------------------------------------------------------------
$ws.Range("A19").Value2 = "Critério:"

# ---------------------------------------------------------------------------
# Row 20: label -> "Norma de recuperação:" (B/C unchanged)
# ---------------------------------------------------------------------------
$ws.Range("A20").Value2 = "Norma de recuperação:"

# ---------------------------------------------------------------------------
# Row 21: label -> "Bibliografia:" (B/C unchanged), height grows to 120
# ---------------------------------------------------------------------------
$ws.Range("A21").Value2 = "Bibliografia:"
$ws.Rows.Item(21).RowHeight = 120

# ---------------------------------------------------------------------------
# Row 22: label -> "Requisitos:" only (B/C stay empty, height already reset)
# ---------------------------------------------------------------------------
$ws.Range("A22").Value2 = "Requisitos:"

# ---------------------------------------------------------------------------
# Row 23: A empty, gains the first "Requisito fraco" line in B/C (new cells)
# ---------------------------------------------------------------------------
Set-NewCell "B23" "B24" "LOB1206 -  Solos I  (Requisito fraco)`n"
Set-NewCell "C23" "C24" "LOB1206 -  Solos I  (Requisito fraco)`n"
$ws.Rows.Item(23).RowHeight = 30

# ---------------------------------------------------------------------------
# Row 24: B/C becomes the second "Requisito fraco" line
# ---------------------------------------------------------------------------
$ws.Range("B24").Value2 = "LOQ4233 -  Gestão de Negócios  (Requisito fraco)`n"
$ws.Range("C24").Value2 = "LOQ4233 -  Gestão de Negócios  (Requisito fraco)`n"

# ---------------------------------------------------------------------------
# Row 25 no longer exists in the target sheet (dimension shrinks to C24)
# ---------------------------------------------------------------------------
$ws.Rows.Item(25).Delete()
